$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.855.99"
$ws.Range("E2").Value = "  +3.32%  "
$ws.Range("D3").Value = "'1.678.84"
$ws.Range("E3").Value = "  +3.19%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'219.46"
$ws.Range("E5").Value = "  +2.35%  "
$ws.Range("D6").Value = "'0.535"
$ws.Range("E6").Value = "  +3.20%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "'29.00"
$ws.Range("E8").Value = "  +2.05%  "
$ws.Range("E9").Value = "  +2.87%  "
$ws.Range("E10").Value = "  +5.66%  "
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("D12").Value = "'1.921.26"
$ws.Range("D13").Value = "'1.673.71"
$ws.Range("E13").Value = "  +2.91%  "
$ws.Range("E14").Value = "  +7.04%  "
$ws.Range("D15").Value = "'10.05"
$ws.Range("E15").Value = "  +9.18%  "
$ws.Range("D16").Value = "'4.06"
$ws.Range("E16").Value = "  +6.19%  "
$ws.Range("D17").Value = "'30.846.55"
$ws.Range("E17").Value = "  +3.24%  "
$ws.Range("E18").Value = "  +2.48%  "
$ws.Range("D19").Value = "'243.09"
$ws.Range("E19").Value = "  +1.41%  "
$ws.Range("E20").Value = "  +2.58%  "
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("E22").Value = "  +3.01%  "
$ws.Range("E23").Value = "  +2.06%  "
$ws.Range("D24").Value = "'2.17"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("D25").Value = "'159.04"
$ws.Range("E25").Value = "  +0.79%  "
$ws.Range("D26").Value = "'15.82"
$ws.Range("E26").Value = "  +2.55%  "
$ws.Range("E27").Value = "  +2.76%  "
$ws.Range("E28").Value = "  +2.06%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("E30").Value = "  +1.17%  "
$ws.Range("E31").Value = "  +3.63%  "
$ws.Range("E32").Value = "  +3.19%  "
$ws.Range("D33").Value = "'1.517.86"
$ws.Range("E33").Value = "  +6.74%  "
$ws.Range("E34").Value = "  +4.29%  "
$ws.Range("E35").Value = "  +4.54%  "
$ws.Range("D36").Value = "'83.97"
$ws.Range("E36").Value = "  +12.70%  "
$ws.Range("E37").Value = "  +0.51%  "
$ws.Range("E38").Value = "  +9.40%  "
$ws.Range("D39").Value = "'0.0179"
$ws.Range("E39").Value = "  +5.07%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("E41").Value = "  -3.64%  "
$ws.Range("E42").Value = "  +3.06%  "
$ws.Range("E43").Value = "  +1.54%  "
$ws.Range("E44").Value = "  +0.42%  "
$ws.Range("E45").Value = "  +2.36%  "
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("D47").Value = "'5.58"
$ws.Range("E47").Value = "  +4.85%  "
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").Value = "'50.68"
$ws.Range("E48").Value = "  +4.32%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "'1.812.07"
$ws.Range("E49").Value = "  +2.52%  "
$ws.Range("E50").Value = "  +6.62%  "
$ws.Range("D51").Value = "'92.87"
$ws.Range("E51").Value = "  +2.19%  "
